$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "color" column
$ws.Range("E1").Value = "color "

# Add color value for existing row (11N60530 / ameer / kia / 2015)
$ws.Range("E2").Value = "red"

# Add a brand new row of car data
$ws.Range("A3").Value = "17A25060"
$ws.Range("B3").Value = "mohamed "
$ws.Range("C3").Value = "reno"
$ws.Range("D3").Value = 2012
$ws.Range("E3").Value = "gray"

# Move selection to mirror the authored state
$ws.Range("E5").Select()
